$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total)
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "32 / 112"
